function Set-TextCell($range, $value) {
    # Force the write to land as text, even when the string looks numeric
    # (e.g. "88.10", "0.514"), then drop the temporary Text number-format so
    # the cell keeps its original (unstyled) appearance.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws.Range("D2") "40.159.88"
Set-TextCell $ws.Range("E2") "  +0.23%  "

Set-TextCell $ws.Range("D3") "2.235.47"
Set-TextCell $ws.Range("E3") "  +0.85%  "

Set-TextCell $ws.Range("E4") "  -0.11%  "

Set-TextCell $ws.Range("D5") "295.25"
Set-TextCell $ws.Range("E5") "  +1.50%  "

Set-TextCell $ws.Range("D6") "88.10"
Set-TextCell $ws.Range("E6") "  -0.37%  "

Set-TextCell $ws.Range("D7") "0.514"
Set-TextCell $ws.Range("E7") "  -0.35%  "

Set-TextCell $ws.Range("E8") "  -0.03%  "

Set-TextCell $ws.Range("D9") "0.470"
Set-TextCell $ws.Range("E9") "  -0.34%  "

Set-TextCell $ws.Range("D10") "30.71"
Set-TextCell $ws.Range("E10") "  -0.43%  "

Set-TextCell $ws.Range("D11") "51.00"
Set-TextCell $ws.Range("E11") "  +6.73%  "

Set-TextCell $ws.Range("D12") "0.0783"
Set-TextCell $ws.Range("E12") "  -0.24%  "

Set-TextCell $ws.Range("E13") "  +2.92%  "

Set-TextCell $ws.Range("D14") "6.47"
Set-TextCell $ws.Range("E14") "  +0.31%  "

Set-TextCell $ws.Range("D15") "2.577.19"
Set-TextCell $ws.Range("E15") "  +0.71%  "

Set-TextCell $ws.Range("D16") "13.91"
Set-TextCell $ws.Range("E16") "  -0.82%  "

Set-TextCell $ws.Range("D17") "2.276.18"
Set-TextCell $ws.Range("E17") "  +3.04%  "

Set-TextCell $ws.Range("D18") "0.738"
Set-TextCell $ws.Range("E18") "  +1.03%  "

Set-TextCell $ws.Range("D19") "40.072.39"
Set-TextCell $ws.Range("E19") "  +0.22%  "

Set-TextCell $ws.Range("D20") "0.0₃0890"
Set-TextCell $ws.Range("E20") "  +0.33%  "

Set-TextCell $ws.Range("D21") "11.33"
Set-TextCell $ws.Range("E21") "  -3.58%  "

Set-TextCell $ws.Range("D22") "5.81"
Set-TextCell $ws.Range("E22") "  -0.19%  "

Set-TextCell $ws.Range("D23") "65.95"
Set-TextCell $ws.Range("E23") "  +0.44%  "

Set-TextCell $ws.Range("D24") "237.42"
Set-TextCell $ws.Range("E24") "  +0.28%  "

Set-TextCell $ws.Range("E25") "  +0.02%  "

Set-TextCell $ws.Range("D26") "2.48"
Set-TextCell $ws.Range("E26") "  +0.37%  "

Set-TextCell $ws.Range("D27") "1.84"
Set-TextCell $ws.Range("E27") "  -0.62%  "

Set-TextCell $ws.Range("B28") "EthereumClassic"
Set-TextCell $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D28") "23.37"
Set-TextCell $ws.Range("E28") "  +3.33%  "

Set-TextCell $ws.Range("B29") "Toncoin"
Set-TextCell $ws.Range("C29") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D29") "2.15"
Set-TextCell $ws.Range("E29") "  -1.92%  "

Set-TextCell $ws.Range("D30") "9.34"
Set-TextCell $ws.Range("E30") "  +0.77%  "

Set-TextCell $ws.Range("D31") "158.11"
Set-TextCell $ws.Range("E31") "  +3.04%  "

Set-TextCell $ws.Range("D32") "31.95"
Set-TextCell $ws.Range("E32") "  -1.04%  "

Set-TextCell $ws.Range("E33") "  -0.02%  "

Set-TextCell $ws.Range("D34") "5.00"
Set-TextCell $ws.Range("E34") "  +0.53%  "

Set-TextCell $ws.Range("D35") "3.06"
Set-TextCell $ws.Range("E35") "  +7.45%  "

Set-TextCell $ws.Range("D36") "0.0717"
Set-TextCell $ws.Range("E36") "  -0.38%  "

Set-TextCell $ws.Range("E37") "  -3.39%  "

Set-TextCell $ws.Range("E38") "  +0.96%  "

Set-TextCell $ws.Range("B39") "ARBITRUM"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D39") "1.77"
Set-TextCell $ws.Range("E39") "  +3.42%  "

Set-TextCell $ws.Range("B40") "Kaspa"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D40") "0.0997"
Set-TextCell $ws.Range("E40") "  -0.38%  "

Set-TextCell $ws.Range("D41") "15.59"
Set-TextCell $ws.Range("E41") "  -3.87%  "

Set-TextCell $ws.Range("D42") "2.096.09"
Set-TextCell $ws.Range("E42") "  +0.14%  "

Set-TextCell $ws.Range("D43") "3.73"
Set-TextCell $ws.Range("E43") "  -2.89%  "

Set-TextCell $ws.Range("D44") "18.78"
Set-TextCell $ws.Range("E44") "  +5.86%  "

Set-TextCell $ws.Range("D45") "10.16"
Set-TextCell $ws.Range("E45") "  +2.28%  "

Set-TextCell $ws.Range("D46") "0.0271"
Set-TextCell $ws.Range("E46") "  +0.54%  "

Set-TextCell $ws.Range("B47") "NEARProtocol"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D47") "2.74"
Set-TextCell $ws.Range("E47") "  +2.74%  "

Set-TextCell $ws.Range("B48") "ApeXProtocol"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws.Range("D48") "1.95"
Set-TextCell $ws.Range("E48") "  -10.25%  "

Set-TextCell $ws.Range("D49") "2.448.56"
Set-TextCell $ws.Range("E49") "  +0.68%  "

Set-TextCell $ws.Range("D50") "1.49"
Set-TextCell $ws.Range("E50") "  +2.72%  "

Set-TextCell $ws.Range("E51") "  +3.44%  "

